$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $true, $true, $false, $false, $false,
                                   $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Host "WARNING: replace failed for: $find"
    }
}

# --- Title ---
Replace-Text "Unveiling Nature's Symphony of Resilience" "Mathematics: Unlocking the Language of the Universe"

# --- Author name paragraph: "Dr. Catherine Rose" -> "Benjamin James" ---
Replace-Text "Dr. Catherine Rose" "Benjamin James"

# --- Email paragraph ---
Replace-Text "catherine" "benjamin"
Replace-Text "rose@abcademica" "james@educationhub"
Replace-Text "edu" "org"

Write-Host "Simple replacements done"

# --- Big body paragraph (paragraph 5) ---
# Block 1: first three sentences collapse to two sentences.
Replace-Text "In the vast tapestry of life, resilience emerges as a harmonious symphony of adaptations, an orchestra of interconnected elements performing a concerto of survival and growth" "From the intricate patterns of nature to the algorithms that power our digital age, mathematics is a guiding force in the world around us"

Replace-Text " From the intricate mechanisms of cellular regeneration to the remarkable plasticity of ecosystems, the natural world abounds with inspiring stories of resilience that teach us about the astounding capabilities of life to endure, adapt, and even thrive in the face of adversity. As we delve into the depths of this phenomenon, we'll explore resilience across diverse domains, unraveling the secrets of nature's resilience and uncovering profound insights into our own capacity for resilience and adaptation in these challenging times" " It is a language, a tool, and a way of thinking that enables us to understand and navigate the complexities of our universe"

Write-Host "Block1 done"

# Block 2: sentence after first double-break.
Replace-Text "Resilience, as it manifests in nature, is not merely a passive endurance but an active engagement with the challenges posed by the environment" "In the realm of numbers, equations, and geometric shapes, mathematics offers a glimpse into the fundamental order that governs the physical world"

Replace-Text " Organisms, from bacteria to towering trees, possess remarkable abilities to sense and respond to changing conditions, modifying their behaviors, and physiological processes to optimize their chances of survival" " From the precise calculations of physics to the intricate algorithms of computer science, mathematics provides the foundation for many of our modern-day technologies"

Replace-Text " From the delicate dance of bacteria in response to antibiotics to the strategic shedding of leaves by drought-stricken plants, nature showcases a myriad of resilience strategies that are both elegant and effective. This adaptability, a hallmark of life's resilience, provides a glimpse into the extraordinary capacity of living systems to reshape themselves and thrive despite obstacles" " It is a language that allows us to communicate ideas, solve problems, and make predictions about the world around us"

Write-Host "Block2 done"

# Block 3: sentence after second double-break.
Replace-Text "The resilience of ecosystems, intricate webs of interconnected organisms, is a testament to the collective power of diversity" "Furthermore, mathematics is more than just a collection of abstract concepts; it is a tool that empowers us to analyze data, model complex systems, and make informed decisions"

Replace-Text " In these natural communities, resilience resides not in individual organisms but in the collective synergy of species, each contributing its unique strengths and adaptations" " From the financial markets to the medical field, mathematics plays a vital role in shaping our understanding of the world and in driving progress across various disciplines"

Write-Host "Block3 done"

# Block 4: remaining tail of paragraph replaced with new "Body:" section plus
# three new sub-sections (history, branches, everyday-life), using ^l for new
# line breaks (Word COM wildcard code for a manual line break).
$newTail = "^l^lBody:^lMathematics has a rich history that spans across cultures and civilizations. " + `
    "From the ancient Babylonians and Egyptians to the groundbreaking work of mathematicians like Pythagoras, Euclid, and Newton, the study of mathematics has continuously evolved and expanded our understanding of the universe." + `
    "^l^lIn its many branches, from algebra to calculus to statistics, mathematics offers a diverse range of concepts and techniques that can be applied to a wide variety of fields. " + `
    "From the movement of celestial bodies to the flow of fluids, mathematics helps us to understand and predict the behavior of the natural world. " + `
    "In the realm of engineering, mathematics enables us to design structures, machines, and systems that are both efficient and reliable." + `
    "^l^lMoreover, mathematics is not just a subject confined to textbooks and classrooms; it is an integral part of our everyday lives. " + `
    "From the time we tell time to the way we navigate using maps, mathematics is woven into the fabric of our daily existence. " + `
    "It helps us make sense of the world around us, from the patterns in nature to the complexities of human interactions."

Replace-Text " As ecosystems face pressures from climate change, habitat loss, and pollution, the resilience of these complex systems is put to the test, demonstrating the critical importance of preserving biodiversity and promoting ecological balance. Exploring the resilience of ecosystems, we gain valuable insights into the interconnectedness of life and the delicate equilibrium upon which our planet's health depends." $newTail

Write-Host "Block4 done"

# --- Summary paragraph (paragraph 7) ---
Replace-Text "Resilience, pervasive in nature, is a symphony of adaptations and responses to environmental challenges" "Mathematics, with its universal language, empowers us to understand the fundamental order of the universe, solve complex problems, and make informed decisions"

Replace-Text " From cellular regeneration to ecosystem resilience, nature's strategies inspire us to understand our capacity for resilience and adaptation" " Its rich history and diverse branches offer a tapestry of concepts and techniques that have revolutionized our understanding of the world across numerous disciplines"

Replace-Text " Resilience is not mere endurance but an active engagement with adversity, showcasing the plasticity and adaptability of life. Unveiling nature's resilience offers profound insights into the interconnectedness of life, biodiversity, and the delicate equilibrium of our planet's health" " From the intricate workings of nature to the complexities of human society, mathematics continues to be a guiding force in our quest for knowledge and progress"

Write-Host "Summary block done"

# --- Trailing empty paragraph added at the very end of the document body ---
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

Write-Host "Trailing paragraph added"





